$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (DataPeriodName) written first, in the same order the ---
# --- original author typed them, so new shared-string entries land in --
# --- the same slots as the authored workbook. -------------------------
$ws.Range("A1").Value = "DataPeriodName"
$ws.Range("A5").Value = "04. Предыдущий месяц (но не ранее января)"
$ws.Range("A3").Value = "02. С начала года до предыдущего месяца (но не ранее января)"
$ws.Range("A6").Value = "05. Первый квартал (январь - март)"
$ws.Range("A7").Value = "06. Второй квартал (апрель - июнь)"
$ws.Range("A9").Value = "08. Четвёртый квартал (октябрь - декабрь)"
$ws.Range("A10").Value = "08. Первое полугодие (январь - июнь)"
$ws.Range("A11").Value = "09. Второе полугодие (июль - декабрь)"
$ws.Range("A12").Value = "10. 1й, 2й и 3й кварталы (январь - сентябрь)"
$ws.Range("A8").Value = "07. Третий квартал (июль - сентябрь)"
$ws.Range("A2").Value = "01. С начала года до текущего месяца"
$ws.Range("A4").Value = "03. Текущий месяц"
$ws.Range("A13").Value = "11. Весь год (12 месяцев)"

# --- Header row (B1/C1 reuse existing shared strings) -------------------
$ws.Range("B1").Value = "StartMonth"
$ws.Range("C1").Value = "EndMonth"

# --- Remaining columns (StartMonth / EndMonth), all reuse existing ------
# --- shared-string entries so they don't disturb the table order. -------
$ws.Range("B2").Value = "01"
$ws.Range("C2").Value = "this"

$ws.Range("B3").Value = "01"
$ws.Range("C3").Value = "this-1"

$ws.Range("B4").Value = "this"
$ws.Range("C4").Value = "this"

$ws.Range("B5").Value = "this-1"
$ws.Range("C5").Value = "this-1"

$ws.Range("B6").Value = "01"
$ws.Range("C6").Value = "03"

$ws.Range("B7").Value = "04"
$ws.Range("C7").Value = "06"

$ws.Range("B8").Value = "07"
$ws.Range("C8").Value = "09"

$ws.Range("B9").Value = "10"
$ws.Range("C9").Value = "12"

$ws.Range("B10").Value = "01"
$ws.Range("C10").Value = "06"

$ws.Range("B11").Value = "07"
$ws.Range("C11").Value = "12"

$ws.Range("B12").Value = "01"
$ws.Range("C12").Value = "09"

$ws.Range("B13").Value = "01"
$ws.Range("C13").Value = "12"

# --- Column A width (bestFit applied in the source workbook) -----------
# Target stored width is 60.42578125 characters; the ColumnWidth property
# adds the standard 5/6 character padding internally, so subtract it back
# out before assigning.
$ws.Columns.Item(1).ColumnWidth = 60.42578125 - (5 / 6)

# --- Selection -----------------------------------------------------------
$ws.Range("A8").Select()
